# Insert a new weekly price record as row 204 on the single data sheet.
# Every existing row from 204 downward (204..330) shifts down by one
# (204->205, ..., 330->331), growing the used range from A1:R330 to
# A1:R331. The freshly inserted row 204 gets its own data (date, volume,
# min/max/avg price, avg-price-per-unit) while inheriting the constant
# market/category/unit/region metadata already used by this block of
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 204:330 down to 205:331, leaving a blank row 204 behind.
$ws.Rows.Item(204).EntireRow.Insert()

# Populate the newly inserted row 204.
$ws.Range("A204").Value = 10
$ws.Range("B204").Value = "Vega Modelo de Temuco"
$ws.Range("C204").Value = "La Araucanía"
$ws.Range("D204").Value = 44719
$ws.Range("E204").Value = 9
$ws.Range("F204").Value = 100112009
$ws.Range("G204").Value = "Acelga"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 20
$ws.Range("K204").Value = 10000
$ws.Range("L204").Value = 10000
$ws.Range("M204").Value = 10000
$ws.Range("N204").Value = "$/docena de atados (12 kilos)"
$ws.Range("O204").Value = "Provincia de Cautín"
$ws.Range("P204").Value = 833
$ws.Range("Q204").Value = 12
$ws.Range("R204").Value = "Hortaliza"
